# Update "想去人数" (number of people interested) figures for three events
# that appear on both the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 2692   # was 2685
$wsExpo.Range("F6").Value = 255    # was 253
$wsExpo.Range("F7").Value = 384    # was 383

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 2692    # was 2685
$wsAll.Range("F6").Value = 255     # was 253
$wsAll.Range("F9").Value = 384     # was 383
